$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

# The "Metadata" sheet gains a new "Contact" row (second contact: Bob Milius)
# and a new "Jurisdiction" row right after the existing Contact row, which
# pushes Description/Purpose/Copyright/Immutable down by one row.
#
# Insert a blank row at 12, then duplicate the border/fill/wrap formatting
# used by the other data rows (copy row 11's format onto the new row 12).
$meta.Rows.Item(12).Insert()
$meta.Rows.Item(11).Copy()
$meta.Range("A12:B12").PasteSpecial(-4122)

# --- simple value updates ---
$meta.Range("B3").Value = "0.1.7"
$meta.Range("B6").Value = "draft"
$meta.Range("B8").Value = "2024-08-23T10:17:11-05:00"
$meta.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# --- new "Contact" row (Bob Milius) ---
$meta.Range("A11").Value = "Contact"
$meta.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- new "Jurisdiction" row (no display value) ---
$meta.Range("A12").Value = "Jurisdiction"
$meta.Range("B12").Value = ""

# --- "Description" row now carries the same text as "Title" ---
$meta.Range("B13").Value = "Auto Differential panel - Blood (57023-4)"

$wb.Save()
